$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Upper Bounds" column (B) from "Infinite" to "Finite" for the
# rows corresponding to models whose bounds were found to be finite
# during bound propagation / standard form conversion.
$ws.Range("B30").Value = "Finite"
$ws.Range("B31").Value = "Finite"
$ws.Range("B38").Value = "Finite"
$ws.Range("B47").Value = "Finite"
